$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Footer "date" field on the slide master + every slide layout:
#    19/11/1443 -> 20/11/1443
# ---------------------------------------------------------------------------
$oldDate = "19/11/1443"
$newDate = "20/11/1443"

$m = $p.SlideMaster

for ($i = 1; $i -le $m.Shapes.Count; $i++) {
    $sh = $m.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

for ($li = 1; $li -le $m.CustomLayouts.Count; $li++) {
    $lyt = $m.CustomLayouts.Item($li)
    for ($i = 1; $i -le $lyt.Shapes.Count; $i++) {
        $sh = $lyt.Shapes.Item($i)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Slide 1 title: fix casing of "sql server databases" -> "SQL Server Databases"
# ---------------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
for ($i = 1; $i -le $s1.Shapes.Count; $i++) {
    $sh = $s1.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.TextRange.Text -eq "CI/CD for sql server databases in Azure DevOps") {
            $sh.TextFrame.TextRange.Text = "CI/CD for SQL Server Databases in Azure DevOps"
        }
    }
}

# ---------------------------------------------------------------------------
# 3) Slide 3 headline: "What Database DevOps is Important?"
#    -> "Why Database DevOps is Important?"
#    ("What " is replaced by "Why " as its own run, leaving the rest intact)
# ---------------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
for ($i = 1; $i -le $s3.Shapes.Count; $i++) {
    $sh = $s3.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        $tr = $sh.TextFrame.TextRange
        if ($tr.Text -eq "What Database DevOps is Important?") {
            $lead = $tr.Characters(1, 5)
            $lead.Text = "Why "
        }
    }
}
